$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the empty "Design Method" sub-header cells in column B for the rows
# that only have a text label in column C (the n=... summary rows). Doing a
# full Clear (not just ClearContents) drops the now-pointless styled-but-
# empty cell entirely, matching the cleaned-up sheet.
$ws.Range("B5").Clear()
$ws.Range("B9").Clear()
$ws.Range("B13").Clear()
$ws.Range("B17").Clear()

# Narrow column A - it no longer needs to fit the old wide label, so size it
# down to the best-fit width for its contents.
$ws.Columns("A").ColumnWidth = 19.166666666666668

# Leave the cursor parked on I12, where the analysis was last being checked.
$ws.Range("I12").Select()
